$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'67.342.48"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +0.73%  "

$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'3.113.07"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("E4").Value = "  -0.08%  "

$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'575.99"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.15%  "

$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'178.24"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +6.30%  "

$ws.Range("E7").Value = "  -0.03%  "

$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'3.109.37"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +1.30%  "

$ws.Range("E9").Value = "  +0.78%  "

$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'6.54"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +2.48%  "

$ws.Range("E11").Value = "  +1.44%  "

$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.467"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -0.56%  "

$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +0.22%  "

$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'36.47"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +1.45%  "

$ws.Range("E15").Value = "  +0.88%  "

$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'3.633.71"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +1.36%  "

$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'67.312.58"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("E18").Value = "  +0.37%  "

$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'3.112.68"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +1.36%  "

$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'16.46"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  -2.68%  "

$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'486.48"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -0.20%  "

$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'0.690"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("E23").Value = "  +0.11%  "

$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'83.69"
$ws.Range("D24").Style = $style

$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'12.74"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -0.56%  "

$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'2.28"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +2.94%  "

$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'10.48"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +1.94%  "

$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -0.04%  "

$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'7.97"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  +1.84%  "

$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'2.32"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +1.87%  "

$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'2.62"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +0.14%  "

$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'28.11"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +1.70%  "

$ws.Range("E33").Value = "  +0.24%  "

$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.0₃0940"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +3.37%  "

$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -0.06%  "

$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'48.08"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +3.86%  "

$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.951"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("E38").Value = "  -0.92%  "

$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.319"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +5.69%  "

$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'49.21"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -0.05%  "

$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'2.02"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +1.58%  "

$ws.Range("E42").Value = "  +0.25%  "

$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'8.31"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -0.33%  "

$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'2.69"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +8.42%  "

$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'2.789.38"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +1.09%  "

$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'373.45"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +0.82%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'0.0346"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'26.71"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +9.30%  "

$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'135.60"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -0.34%  "

$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'2.37"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +10.34%  "
